$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after row 317 (before old row 318),
# which pushes the old rows 318..342 down to 320..344.
$ws.Rows.Item(318).Insert()
$ws.Rows.Item(318).Insert()

# Populate common/fixed columns for the two new rows, mirroring row 317.
$rows = @(318, 319)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 1).Value2 = 10
    $ws.Cells.Item($r, 2).Value2 = "Vega Modelo de Temuco"
    $ws.Cells.Item($r, 3).Value2 = "La Araucanía"
    $ws.Cells.Item($r, 5).Value2 = 9
    $ws.Cells.Item($r, 6).Value2 = 100112017
    $ws.Cells.Item($r, 7).Value2 = "Apio"
    $ws.Cells.Item($r, 8).Value2 = "Americana (o)"
    $ws.Cells.Item($r, 14).Value2 = "`$/docena de matas"
    $ws.Cells.Item($r, 15).Value2 = "Provincia del Elquí"
    $ws.Cells.Item($r, 17).Value2 = 6
    $ws.Cells.Item($r, 18).Value2 = "Hortaliza"
}

# Row 318: new data point
$ws.Cells.Item(318, 4).Value2 = 44783
$ws.Range("D318").NumberFormat = $ws.Range("D320").NumberFormat
$ws.Cells.Item(318, 9).Value2 = "Primera"
$ws.Cells.Item(318, 10).Value2 = 80
$ws.Cells.Item(318, 11).Value2 = 12000
$ws.Cells.Item(318, 12).Value2 = 12000
$ws.Cells.Item(318, 13).Value2 = 12000
$ws.Cells.Item(318, 16).Value2 = 2000

# Row 319: new data point
$ws.Cells.Item(319, 4).Value2 = 44783
$ws.Range("D319").NumberFormat = $ws.Range("D320").NumberFormat
$ws.Cells.Item(319, 9).Value2 = "Segunda"
$ws.Cells.Item(319, 10).Value2 = 30
$ws.Cells.Item(319, 11).Value2 = 10000
$ws.Cells.Item(319, 12).Value2 = 10000
$ws.Cells.Item(319, 13).Value2 = 10000
$ws.Cells.Item(319, 16).Value2 = 1667
